$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOR_Login")
Write-Host $ws.Name
